$d = $word.ActiveDocument

# Locate the paragraph that currently reads "This is ready to submit."
# (Range.Text includes the trailing paragraph-mark character, so match
# with TrimEnd instead of an exact equality check.)
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "This is ready to submit.") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph 'This is ready to submit.'"
}

$r = $target.Range

# Split the single run "This is ready to submit." into two runs:
# "This is ready to submit" and "." - by replacing the paragraph's XML
# contents directly (InsertXML only overwrites what is inside the
# addressed range, so the paragraph's own attributes below are what will
# persist on save).
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
       'w14:paraId="2D337840" w14:textId="4630D8A6" w:rsidR="00C14391" ' +
       'w:rsidRPr="006A6B59" w:rsidRDefault="00C14391" w:rsidP="0023597E">' +
       '<w:r><w:t>This is ready to submit</w:t></w:r>' +
       '<w:r><w:t>.</w:t></w:r>' +
       '</w:p>'
$r.InsertXML($xml)

# Re-fetch the (now split) paragraph and add a new paragraph right after it
# containing "This is a sample project I am about to submit."
$target2 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "This is ready to submit.") {
        $target2 = $p
        break
    }
}

$target2.Range.InsertParagraphAfter()
$newPara = $target2.Next()
$newPara.Range.Text = "This is a sample project I am about to submit."
